# Update column G ("K") values for rows 2-13 as part of regenerating
# save_data to use K (strikeouts) instead of Strike# and recalculated
# std/mean s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 3
    4  = 0
    5  = 4
    6  = 3
    7  = 5
    8  = 1
    9  = 6
    10 = 3
    11 = 2
    12 = 2
    13 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
